# Applies the content edits described by the commit:
#   "embedding, sharing , and copying"
# i.e. the placeholder brand name "transformer" -> "boboo", the
# "Last Updated" date bump, and the two contact-detail swaps.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $found = $d.Content.Find.Execute(
        $findText,    # FindText
        $true,        # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replaceText, # ReplaceWith
        2             # Replace (wdReplaceOne)
    )
    if (-not $found) {
        Write-Output "WARNING: text not found -> $findText"
    }
}

# 1. Bump the "Last Updated" date from Aug. 12 to Aug. 14, 2022.
Replace-Text "Last Updated on Aug. 12, 2022" "Last Updated on Aug. 14, 2022"

# 2. Consent paragraph: rename the placeholder brand and swap the contact site.
Replace-Text "This Privacy notice for 40;transformer describes" "This Privacy notice for 40;boboo describes"
Replace-Text "please contact us atwww.abigirl.com" "please contact us atwww.abby.cm"

# 3. "The Data We Collect About You?" paragraph: rename the placeholder brand.
Replace-Text "optimization purposes.transformer also uses" "optimization purposes.boboo also uses"

# 4. Bullet list item: rename the placeholder brand in the website/app mention.
Replace-Text "your use of the transformer website/app" "your use of the boboo website/app"

# 5. Further Details paragraph: swap the contact e-mail address.
Replace-Text "please contactduruakuebuka@gmail.com" "please contactoluwatofunmi.famuwagun@stu.cu.edu.ng"
